$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the text of the four existing comments before we start shuffling
#    rows around (comments stay anchored to their absolute cell address and
#    are NOT moved automatically by row insert/delete operations).
# ---------------------------------------------------------------------------
$commentCamera       = $ws.Range("B2").Comment.Text()    # on the row being removed
$commentGameObject   = $ws.Range("B5").Comment.Text()    # -> should end up on B3
$commentShaders      = $ws.Range("B13").Comment.Text()   # -> should end up on B11
$commentErrHandling  = $ws.Range("B18").Comment.Text()   # -> should end up on B17

# ---------------------------------------------------------------------------
# 2. Remove the four existing comments now - we'll re-create the three that
#    survive at their new locations further down, and simply drop the one
#    that belonged to the deleted row.
# ---------------------------------------------------------------------------
$ws.Range("B2").Comment.Delete()
$ws.Range("B5").Comment.Delete()
$ws.Range("B13").Comment.Delete()
$ws.Range("B18").Comment.Delete()

# ---------------------------------------------------------------------------
# 3. Remember the "On-screen debug text" row's data (it moves from its
#    current spot to just after the "Compiled shaders" row, with its
#    estimate changed from 14 to 21).
# ---------------------------------------------------------------------------
$debugTextCol1 = $ws.Cells.Item(3, 1).Value()
$debugTextCol2 = $ws.Cells.Item(3, 2).Value()

# ---------------------------------------------------------------------------
# 4. Delete row 2 ("Model Viewer" / "Add camera controls to model viewer").
#    This shifts every following row up by one, so the "On-screen debug
#    text" row (previously row 3) is now row 2.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 5. Delete the "On-screen debug text" row (now row 2) - this removes it
#    from its old position and shifts everything else up by one again.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 6. Insert a new row right after "Compiled shaders" (which is now row 11)
#    and populate it with the "On-screen debug text" task, using the new
#    estimate of 21.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12, 1).Value = $debugTextCol1
$ws.Cells.Item(12, 2).Value = $debugTextCol2
$ws.Cells.Item(12, 3).Value = 21

# ---------------------------------------------------------------------------
# 7. Re-create the three surviving comments at their new cell addresses.
# ---------------------------------------------------------------------------
$ws.Range("B3").AddComment($commentGameObject)
$ws.Range("B11").AddComment($commentShaders)
$ws.Range("B17").AddComment($commentErrHandling)

# ---------------------------------------------------------------------------
# 8. Fix up the sheet view selection to match the new target state.
# ---------------------------------------------------------------------------
$ws.Range("C12").Select()
